$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "need_to_buy" data table: shift the day index forward by 10
# and refresh the TISG/fcs/buy/MYDIR/need_to_buy values for rows 2-15.
$data = @{
    2 = @(46020, 10343.2351019585, 10490.9753871234, 19176.26, 7674.53599709828, -42.114525657429)
    3 = @(46021, 10608.1372457288, 10763.5256434299, 11232.26, 8318.54290949053, 327.075356371683)
    4 = @(46022, 11099.979711269, 11540.9719032307, 11232.26, 8700.11370218618, 375.367733559038)
    5 = @(46023, 5204.49985291439, 9274.59574665926, 12075.86, 8884.36483926792, 253.462524413632)
    6 = @(46024, 11665.8765888327, 11163.8645522205, 12075.86, 9062.51869761636, 339.605135409868)
    7 = @(46025, 5070.88502434833, 7901.5423347254, 12075.86, 8513.19344548277, 180.786490842007)
    8 = @(46026, 4937.82896317932, 7277.16692251248, 12075.86, 8502.56105557318, 154.327832420236)
    9 = @(46027, 13230.991612671, 11225.8358864622, 12075.86, 9162.0747185107, 346.335441873872)
    10 = @(46028, 5297.15477007252, 7545.42959573316, 12075.86, 8740.766434472, 175.430667925215)
    11 = @(46029, 13230.991612671, 12019.8179604736, 12075.86, 9162.0747185107, 379.418028291012)
    12 = @(46030, 13230.991612671, 12258.1105489932, 12075.86, 9162.0747185107, 389.346886145996)
    13 = @(46031, 13230.991612671, 11616.8392614773, 12075.86, 9162.0747185107, 362.627249166165)
    14 = @(46032, 5463.01697963343, 8286.51203711876, 12075.86, 8753.045892182, 206.820747054199)
    15 = @(46033, 5297.15477007252, 8267.50461035088, 12075.86, 8740.766434472, 205.51712686762)
}

# Rows 3 and 4 keep the same "buy" (column D) value as before, so that
# cell is intentionally left untouched for those two rows.
$unchangedD = @(3, 4)

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 1).Value = $vals[0]
    $ws.Cells.Item($r, 2).Value = $vals[1]
    $ws.Cells.Item($r, 3).Value = $vals[2]
    if (-not ($unchangedD -contains $r)) {
        $ws.Cells.Item($r, 4).Value = $vals[3]
    }
    $ws.Cells.Item($r, 5).Value = $vals[4]
    $ws.Cells.Item($r, 6).Value = $vals[5]
}
